{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Identify and delete the paragraph immediately following the\n// \"LOB1012: Estat\u00edstica (Requisito fraco)\" requirement line, together with\n// the \"Ver no Jupiter Salvar em pdf Salvar em docx\" line and the\n// \"\u00a9 2020 . Contact: luizeleno@usp.br. ...\" footer line that follow it.\nconst items = paragraphs.items;\nconst toDelete = [];\nfor (let i = 0; i < items.length; i++) {\n  const text = items[i].text;\n  if (\n    text.indexOf(\"Ver no Jupiter Salvar em pdf Salvar em docx\") !== -1 ||\n    text.indexOf(\"Contact: luizeleno@usp.br\") !== -1\n  ) {\n    toDelete.push(items[i]);\n    // Also remove the blank paragraph immediately preceding the\n    // \"Ver no Jupiter...\" paragraph (i.e. the one right after LOB1012).\n    if (text.indexOf(\"Ver no Jupiter Salvar em pdf Salvar em docx\") !== -1 && i > 0) {\n      const prevText = items[i - 1].text;\n      if (prevText.trim() === \"\") {\n        toDelete.push(items[i - 1]);\n      }\n    }\n  }\n}\n\nfor (const p of toDelete) {\n  p.delete();\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the \"Ver no Jupiter Salvar em pdf Salvar em docx\" paragraph and the\n# \"\u00a9 2020 . Contact: luizeleno@usp.br. ...\" paragraph that follows it, plus\n# the blank paragraph right before them (which sits right after the\n# \"LOB1012: Estat\u00edstica (Requisito fraco)\" requirement line), and remove all\n# three paragraphs.\n$count = $d.Paragraphs.Count\n$indices = @()\nfor ($i = 1; $i -le $count; $i++) {\n  $t = $d.Paragraphs.Item($i).Range.Text\n  if ($t -like \"*Ver no Jupiter Salvar em pdf Salvar em docx*\") {\n    $indices += ($i - 1)\n    $indices += $i\n  }\n  if ($t -like \"*Contact: luizeleno@usp.br*\") {\n    $indices += $i\n  }\n}\n\n$indices = $indices | Sort-Object -Unique -Descending\n\nforeach ($idx in $indices) {\n  $d.Paragraphs.Item($idx).Range.Delete()\n}\n"}
